$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells where the numeric-looking text must be preserved exactly
# (Excel would otherwise normalize/round the number and lose the original
# formatting), so force a Text number format before assigning the value.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.689.36'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '1.879.89'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("E4").Value = '  -0.58%  '
$ws.Range("D5").Value = '331.53'
$ws.Range("E5").Value = '  +2.53%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("D7").Value = '0.4724'
$ws.Range("E7").Value = '  +5.12%  '
$ws.Range("E8").Value = '  +2.92%  '
$ws.Range("D9").Value = '47.99'
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("D10").Value = '0.08054'
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").Value = '21.85'
$ws.Range("E12").Value = '  +1.70%  '
$ws.Range("D13").Value = '1.879.14'
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").Value = '5.967'
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("D17").Value = '87.12'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").Value = '0.00001044'
$ws.Range("E18").Value = '  +0.65%  '
$ws.Range("D19").Value = '0.06620'
$ws.Range("E19").Value = '  +1.12%  '
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("D22").Value = '27.689.73'
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("D23").Value = '5.515'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = '11.03'
$ws.Range("E24").Value = '  +1.49%  '
$ws.Range("D25").Value = '2.295'
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").Value = '2.099.20'
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").Value = '156.31'
$ws.Range("E27").Value = '  +3.33%  '
$ws.Range("D28").Value = '20.26'
$ws.Range("E28").Value = '  +4.21%  '
$ws.Range("E29").Value = '  +3.05%  '
$ws.Range("D30").Value = '5.609'
$ws.Range("E30").Value = '  +1.44%  '
$ws.Range("D31").Value = '122.45'
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("D32").Value = '0.9723'
$ws.Range("E32").Value = '  +4.77%  '
$ws.Range("D33").Value = '0.09576'
$ws.Range("E33").Value = '  +1.88%  '
$ws.Range("D34").Value = '1.450'
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("D35").Value = '3.628'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  +0.79%  '
$ws.Range("D37").Value = '0.06121'
$ws.Range("E37").Value = '  +2.37%  '
$ws.Range("E38").Value = '  +1.09%  '
$ws.Range("D39").Value = '1.235'
$ws.Range("E39").Value = '  +1.09%  '
$ws.Range("D40").Value = '8.154'
$ws.Range("E40").Value = '  -2.88%  '
$ws.Range("E41").Value = '  +1.21%  '
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("D43").Value = '0.1905'
$ws.Range("E43").Value = '  +2.65%  '
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").Value = '0.5704'
$ws.Range("E45").Value = '  +0.90%  '
$ws.Range("D46").Value = '1.245'
$ws.Range("E46").Value = '  -2.54%  '
$ws.Range("D47").Value = '12.23'
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").Value = '3.409'
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("D49").Value = '1.933'
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("E50").Value = '  +9.35%  '
$ws.Range("D51").Value = '0.06821'
$ws.Range("E51").Value = '  -0.67%  '
